# Apply updated crypto market data (price/volume refresh) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    # Force text interpretation so numeric-looking strings (e.g. "0.7124")
    # are not silently converted to floating point numbers by Excel,
    # matching the inlineStr text cells used in the source workbook.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '29.194.35'
Set-TextCell 2 5 '  -1.18%  '

# Row 3
Set-TextCell 3 4 '1.862.53'
Set-TextCell 3 5 '  -0.97%  '

# Row 4
Set-TextCell 4 5 '  +0.28%  '

# Row 5
Set-TextCell 5 4 '0.7124'
Set-TextCell 5 5 '  -1.61%  '

# Row 6
Set-TextCell 6 4 '240.49'
Set-TextCell 6 5 '  +0.50%  '

# Row 7
Set-TextCell 7 5 '  +0.23%  '

# Row 8: 'Cardano' -> 'Dogecoin'
Set-TextCell 8 2 'Dogecoin'
Set-TextCell 8 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 8 4 '0.07699'
Set-TextCell 8 5 '  -2.39%  '

# Row 9: 'Dogecoin' -> 'Cardano'
Set-TextCell 9 2 'Cardano'
Set-TextCell 9 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 9 4 '0.3073'
Set-TextCell 9 5 '  -0.48%  '

# Row 10
Set-TextCell 10 4 '24.90'
Set-TextCell 10 5 '  -1.85%  '

# Row 11
Set-TextCell 11 4 '0.08263'
Set-TextCell 11 5 '  +0.63%  '

# Row 12
Set-TextCell 12 4 '1.855.58'
Set-TextCell 12 5 '  -0.75%  '

# Row 13
Set-TextCell 13 4 '0.7166'
Set-TextCell 13 5 '  -1.40%  '

# Row 14
Set-TextCell 14 4 '5.212'
Set-TextCell 14 5 '  -1.15%  '

# Row 15
Set-TextCell 15 4 '90.18'
Set-TextCell 15 5 '  +0.58%  '

# Row 16
Set-TextCell 16 4 '29.177.62'
Set-TextCell 16 5 '  -1.03%  '

# Row 17
Set-TextCell 17 4 '5.857'
Set-TextCell 17 5 '  +0.06%  '

# Row 18
Set-TextCell 18 4 '243.31'
Set-TextCell 18 5 '  +0.49%  '

# Row 19
Set-TextCell 19 4 '0.000007796'
Set-TextCell 19 5 '  -1.00%  '

# Row 20
Set-TextCell 20 4 '13.14'
Set-TextCell 20 5 '  -1.87%  '

# Row 21
Set-TextCell 21 4 '2.109.17'
Set-TextCell 21 5 '  -0.07%  '

# Row 22
Set-TextCell 22 5 '  +0.14%  '

# Row 23
Set-TextCell 23 5 '  +2.51%  '

# Row 24
Set-TextCell 24 5 '  +0.29%  '

# Row 25
Set-TextCell 25 4 '0.1577'
Set-TextCell 25 5 '  +6.03%  '

# Row 26
Set-TextCell 26 4 '162.23'
Set-TextCell 26 5 '  -0.40%  '

# Row 27
Set-TextCell 27 4 '8.902'
Set-TextCell 27 5 '  -1.08%  '

# Row 28
Set-TextCell 28 4 '18.20'
Set-TextCell 28 5 '  -0.47%  '

# Row 29
Set-TextCell 29 5 '  -2.91%  '

# Row 30
Set-TextCell 30 4 '1.496'
Set-TextCell 30 5 '  +0.89%  '

# Row 31
Set-TextCell 31 4 '4.347'
Set-TextCell 31 5 '  -0.20%  '

# Row 32
Set-TextCell 32 4 '4.087'
Set-TextCell 32 5 '  -0.52%  '

# Row 33
Set-TextCell 33 4 '0.05178'
Set-TextCell 33 5 '  -1.33%  '

# Row 34
Set-TextCell 34 4 '1.901'
Set-TextCell 34 5 '  -2.73%  '

# Row 35
Set-TextCell 35 5 '  -2.01%  '

# Row 36
Set-TextCell 36 4 '0.7266'
Set-TextCell 36 5 '  +1.06%  '

# Row 37
Set-TextCell 37 4 '2.684'
Set-TextCell 37 5 '  +0.43%  '

# Row 38
Set-TextCell 38 4 '0.01847'
Set-TextCell 38 5 '  -0.76%  '

# Row 39
Set-TextCell 39 4 '2.686'
Set-TextCell 39 5 '  -0.89%  '

# Row 40
Set-TextCell 40 4 '1.141.82'
Set-TextCell 40 5 '  -3.14%  '

# Row 41
Set-TextCell 41 4 '0.8994'
Set-TextCell 41 5 '  -1.40%  '

# Row 42
Set-TextCell 42 4 '6.074'
Set-TextCell 42 5 '  +1.35%  '

# Row 43
Set-TextCell 43 4 '72.09'
Set-TextCell 43 5 '  +0.06%  '

# Row 44
Set-TextCell 44 5 '  +0.19%  '

# Row 45
Set-TextCell 45 4 '101.37'
Set-TextCell 45 5 '  -1.08%  '

# Row 46
Set-TextCell 46 4 '0.5266'
Set-TextCell 46 5 '  -1.48%  '

# Row 47
Set-TextCell 47 4 '2.004.97'
Set-TextCell 47 5 '  -0.39%  '

# Row 48
Set-TextCell 48 4 '1.761'
Set-TextCell 48 5 '  -1.01%  '

# Row 49: 'BabyDogeCoin' -> 'EnergySwap'
Set-TextCell 49 2 'EnergySwap'
Set-TextCell 49 3 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 49 4 '9.298'
Set-TextCell 49 5 '  +0.76%  '

# Row 50: 'EnergySwap' -> 'SynthetixNetwork'
Set-TextCell 50 2 'SynthetixNetwork'
Set-TextCell 50 3 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
Set-TextCell 50 4 '2.862'
Set-TextCell 50 5 '  -0.81%  '

# Row 51: 'SynthetixNetwork' -> 'Frax'
Set-TextCell 51 2 'Frax'
Set-TextCell 51 3 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextCell 51 4 '1.000'
Set-TextCell 51 5 '  -0.04%  '
